$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F (shifts the old "test" column from F to G)
$ws.Columns("F").Insert()

# New column F header + fill in the "explicit null cells" row; row 2 (empty
# cells row) and row 3 (empty cells row) are intentionally left blank.
$ws.Range("F1").Value = "default_value"
$ws.Range("F4").Value = "null"

# Match the default column width used for the new default_value column
$ws.Columns("F").ColumnWidth = 13.6

# C2 (unix_format value) loses its explicit number-format style
$ws.Range("C2").Style = "Normal"

# Page orientation now explicit
$ws.PageSetup.Orientation = 1

# Selection now spans the whole used range
$ws.Range("A1:G4").Select() | Out-Null

Write-Output "done"
